$d = $word.ActiveDocument

# Update the iteration "recette" dates: 18 mai -> 17 mai, 8 juin -> 7 juin, 29 juin -> 28 juin
$r = $d.Content
$r.Find.Execute("18 mai, le 8 juin, et le 29 juin", $true, $false, $false, $false, $false,
                 $true, 1, $false, "17 mai, le 7 juin, et le 28 juin", 2)

# Re-anchor the "_GoBack" bookmark right after the last edited digit ("28"),
# matching where Word leaves it after the most recent edit.
$r2 = $d.Content
$r2.Find.Execute("28 juin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($r2.Start + 2, $r2.Start + 2)
$d.Bookmarks.Add("_GoBack", $target)

# The page-count footer on the last section shows the PAGE field's cached value;
# repagination moved it from page 7 to page 5 (NUMPAGES stays at 7).
$sec = $d.Sections.Item($d.Sections.Count)
$footer = $sec.Footers.Item(1)
$footer.Range.Find.Execute("7", $true, $false, $false, $false, $false,
                            $true, 1, $false, "5", 1)
